$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new observation row 20
$ws.Range("B20").Value = 52.76
$ws.Range("C20").Value = 51

# Extend the shared "humidity offset" formula (B - C) down to row 20
$ws.Range("D20").Formula = "=B20-C20"

# Reflect the new active selection
$ws.Range("D20").Select()
